# Generate Report for Archive
#
# Changes:
#   1. Shared string "Ready for handoff" -> "In Translation" everywhere it is
#      used (Status column on the Overview / zh-cn / de-de sheets).
#   2. The Status column's width shrinks to fit the new (shorter) text on
#      each of the three sheets: Overview columns E & F, and column C on
#      both zh-cn and de-de.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- 1. Update the status text -------------------------------------------
$overview.Range("E2:F4").Value = "In Translation"
$zhcn.Range("C2:C4").Value = "In Translation"
$dede.Range("C2:C4").Value = "In Translation"

# --- 2. Resize the affected columns to the new auto-fit width ------------
# Target stored column width (OOXML "characters" units) is 13.4101845877511.
# Excel's ColumnWidth property is expressed in the same "characters" units
# minus a constant 5/6 padding offset that gets re-added when the width is
# persisted, so back that offset out before assigning.
$targetColumnWidth = 13.4101845877511 - (5 / 6)

$overview.Columns.Item(5).ColumnWidth = $targetColumnWidth   # column E
$overview.Columns.Item(6).ColumnWidth = $targetColumnWidth   # column F
$zhcn.Columns.Item(3).ColumnWidth = $targetColumnWidth        # column C
$dede.Columns.Item(3).ColumnWidth = $targetColumnWidth        # column C
